$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the R3 leg row (row 9) values per the calibration fix
$ws.Range("D9").Value = 1300
$ws.Range("G9").Value = 850

# Update the active view/selection to match the saved state
$ws.Range("G11").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
